# Generate Report for Handback
# Adds a new handback record (ef757400-f198-43a2-89dc-3a2f8efce42a.md) as
# row 3 on each of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$newFile        = "ef757400-f198-43a2-89dc-3a2f8efce42a.md"
$newPath        = "e2e\ef757400-f198-43a2-89dc-3a2f8efce42a.md"
$newExt         = ".md"
$status         = "Handed back: in sync with en-US"
$overviewDate   = "2016-09-01 05:05:52"

$zhXlf          = "ef757400-f198-43a2-89dc-3a2f8efce42a.4469a46bbf4f08b79d1fc19f248923d3b56e59ff.zh-cn.xlf"
$zhHandoffDate  = "2016-09-01 05:05:48"
$zhHandbackDate = "2016-09-01 05:06:14"

$deXlf          = "ef757400-f198-43a2-89dc-3a2f8efce42a.4469a46bbf4f08b79d1fc19f248923d3b56e59ff.de-de.xlf"
$deHandoffDate  = "2016-09-01 05:05:52"
$deHandbackDate = "2016-09-01 05:06:22"

$srcRepoUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/12999da2be587f656c4069d9b0a6a61253644f3b/e2e/$newFile"
$zhRepoUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a9843599061cfa0b6f2be6a9434fb3eff65540f3/e2e/$newFile"
$deRepoUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fcf6bc42140ac33ef2646c7ff8aa6f2bf8fe7429/e2e/$newFile"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = $newPath
$wsOverview.Range("C3").Value = $newExt
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $overviewDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $srcRepoUrl, "", "", $newPath) | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()

$wsZh.Range("A3").Value = $newFile
$wsZh.Range("B3").Value = $newExt
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("I3").Value = $newFile
$wsZh.Range("J3").Value = $zhXlf
$wsZh.Range("K3").Value = $zhHandbackDate
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $srcRepoUrl, "", "", $newFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhRepoUrl, "", "", $newFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()

$wsDe.Range("A3").Value = $newFile
$wsDe.Range("B3").Value = $newExt
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("I3").Value = $newFile
$wsDe.Range("J3").Value = $deXlf
$wsDe.Range("K3").Value = $deHandbackDate
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $srcRepoUrl, "", "", $newFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deRepoUrl, "", "", $newFile) | Out-Null

# ---------------------------------------------------------------------
# Update the existing "fa571b25..." record's timestamps/filenames to the
# newly regenerated values (d8cd1e55... files), per the handback rerun.
# ---------------------------------------------------------------------
$oldFile = "fa571b25-7f82-4861-9334-3d0e007235bd.md"
$oldPath = "e2e\fa571b25-7f82-4861-9334-3d0e007235bd.md"
$renFile = "d8cd1e55-3867-4579-8989-057f168cb3f0.md"
$renPath = "e2e\d8cd1e55-3867-4579-8989-057f168cb3f0.md"

$wsOverview.Range("A2").Value = $renFile
$wsOverview.Range("B2").Value = $renPath
$wsOverview.Range("G2").Value = $overviewDate

$wsZh.Range("A2").Value = $renFile
$wsZh.Range("G2").Value = "d8cd1e55-3867-4579-8989-057f168cb3f0.00b8ba9f844c32eb982936f7e9a7a48c97bdda95.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-01 05:05:48"
$wsZh.Range("I2").Value = $renFile
$wsZh.Range("J2").Value = "d8cd1e55-3867-4579-8989-057f168cb3f0.00b8ba9f844c32eb982936f7e9a7a48c97bdda95.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-01 05:06:14"

$wsDe.Range("A2").Value = $renFile
$wsDe.Range("G2").Value = "d8cd1e55-3867-4579-8989-057f168cb3f0.00b8ba9f844c32eb982936f7e9a7a48c97bdda95.de-de.xlf"
$wsDe.Range("H2").Value = $overviewDate
$wsDe.Range("I2").Value = $renFile
$wsDe.Range("J2").Value = "d8cd1e55-3867-4579-8989-057f168cb3f0.00b8ba9f844c32eb982936f7e9a7a48c97bdda95.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-01 05:06:22"

Write-Output "Handback report regenerated: 2 source files."
